# coefficients.xlsx — retext the 0.5 coefficient to 0.25, introduce the new 0.05/0.2
# coefficients, and fill in route.next predictions for ids 4-10 (rows 5-11).
#
# All the numeric-looking coefficients in this sheet are stored as TEXT (shared
# strings), not real numbers (see the still-numeric literal 0 cells in B4:E4 /
# F3:I3 for contrast). Plain `Range.Value = "0.25"` gets auto-coerced back to a
# real number by the engine, so each distinct coefficient is "stamped" once as")
# text via a throwaway Text number format, then fanned out to every target cell
# with Copy / PasteSpecial(xlPasteValues) — a values-only paste carries the
# source's text-ness along without leaving a NumberFormat behind on the copies.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteValues = -4163

# --- Prepare one clean text "stamp" cell per distinct coefficient value ---
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "0.25"
$ws.Range("B3").ClearFormats()
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "0.05"
$ws.Range("F8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2"
$ws.Range("D9").ClearFormats()

# --- Fan each stamp out to every cell that holds that coefficient ---
$ws.Range("B2").Copy()
$ws.Range("F6").PasteSpecial($xlPasteValues)
$ws.Range("G6").PasteSpecial($xlPasteValues)
$ws.Range("H6").PasteSpecial($xlPasteValues)
$ws.Range("I6").PasteSpecial($xlPasteValues)
$ws.Range("B7").PasteSpecial($xlPasteValues)
$ws.Range("C7").PasteSpecial($xlPasteValues)
$ws.Range("D7").PasteSpecial($xlPasteValues)
$ws.Range("E7").PasteSpecial($xlPasteValues)
$ws.Range("B9").PasteSpecial($xlPasteValues)
$ws.Range("C9").PasteSpecial($xlPasteValues)
$ws.Range("H9").PasteSpecial($xlPasteValues)
$ws.Range("I9").PasteSpecial($xlPasteValues)

$ws.Range("B3").Copy()
$ws.Range("C3").PasteSpecial($xlPasteValues)
$ws.Range("D3").PasteSpecial($xlPasteValues)
$ws.Range("E3").PasteSpecial($xlPasteValues)
$ws.Range("F4").PasteSpecial($xlPasteValues)
$ws.Range("G4").PasteSpecial($xlPasteValues)
$ws.Range("H4").PasteSpecial($xlPasteValues)
$ws.Range("I4").PasteSpecial($xlPasteValues)
$ws.Range("B5").PasteSpecial($xlPasteValues)
$ws.Range("C5").PasteSpecial($xlPasteValues)
$ws.Range("D5").PasteSpecial($xlPasteValues)
$ws.Range("E5").PasteSpecial($xlPasteValues)
$ws.Range("F5").PasteSpecial($xlPasteValues)
$ws.Range("G5").PasteSpecial($xlPasteValues)
$ws.Range("H5").PasteSpecial($xlPasteValues)
$ws.Range("I5").PasteSpecial($xlPasteValues)
$ws.Range("B6").PasteSpecial($xlPasteValues)
$ws.Range("C6").PasteSpecial($xlPasteValues)
$ws.Range("D6").PasteSpecial($xlPasteValues)
$ws.Range("E6").PasteSpecial($xlPasteValues)
$ws.Range("F7").PasteSpecial($xlPasteValues)
$ws.Range("G7").PasteSpecial($xlPasteValues)
$ws.Range("H7").PasteSpecial($xlPasteValues)
$ws.Range("I7").PasteSpecial($xlPasteValues)
$ws.Range("B10").PasteSpecial($xlPasteValues)
$ws.Range("C10").PasteSpecial($xlPasteValues)
$ws.Range("D10").PasteSpecial($xlPasteValues)
$ws.Range("E10").PasteSpecial($xlPasteValues)
$ws.Range("F11").PasteSpecial($xlPasteValues)
$ws.Range("G11").PasteSpecial($xlPasteValues)
$ws.Range("H11").PasteSpecial($xlPasteValues)
$ws.Range("I11").PasteSpecial($xlPasteValues)

$ws.Range("F8").Copy()
$ws.Range("B8").PasteSpecial($xlPasteValues)
$ws.Range("C8").PasteSpecial($xlPasteValues)
$ws.Range("D8").PasteSpecial($xlPasteValues)
$ws.Range("E8").PasteSpecial($xlPasteValues)
$ws.Range("G8").PasteSpecial($xlPasteValues)
$ws.Range("H8").PasteSpecial($xlPasteValues)
$ws.Range("I8").PasteSpecial($xlPasteValues)
$ws.Range("F10").PasteSpecial($xlPasteValues)
$ws.Range("G10").PasteSpecial($xlPasteValues)
$ws.Range("H10").PasteSpecial($xlPasteValues)
$ws.Range("I10").PasteSpecial($xlPasteValues)
$ws.Range("B11").PasteSpecial($xlPasteValues)
$ws.Range("C11").PasteSpecial($xlPasteValues)
$ws.Range("D11").PasteSpecial($xlPasteValues)
$ws.Range("E11").PasteSpecial($xlPasteValues)

$ws.Range("D9").Copy()
$ws.Range("E9").PasteSpecial($xlPasteValues)
$ws.Range("F9").PasteSpecial($xlPasteValues)
$ws.Range("G9").PasteSpecial($xlPasteValues)

$excel.CutCopyMode = 0

# Selection moves to B10 in the saved view, matching the authored edit
$ws.Range("B10").Select()
